$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new trade record
$ws.Range("A5").Value = 9582.34
$ws.Range("B5").Value = 9864.4599999999991
$ws.Range("C5").Value = 113.86
$ws.Range("D5").Value = 110.6
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -2.86
$ws.Range("G5").Value = 42607.884212962963
$ws.Range("H5").Value = $false

# Row 6: new trade record
$ws.Range("A6").Value = 9335.1200000000008
$ws.Range("B6").Value = 9582.34
$ws.Range("C6").Value = 110.77
$ws.Range("D6").Value = 107.91
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -2.58
$ws.Range("G6").Value = 42608.616296296299
$ws.Range("H6").Value = $false
